$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the manning (F) column values for rows 2-15 to 0.04
$ws.Range("F2:F15").Value = 0.04

# Select F3:F15 with F3 as the active cell, matching the saved view state
$ws.Activate()
$ws.Range("F3:F15").Select()
$excel.ActiveCell = $ws.Range("F3")
